$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin price (column D) and 1h volume change (column E) values.
# Column D occasionally holds plain-looking numeric strings (e.g. "1.00", "6.40");
# force those to remain plain text (matching the original inlineStr cells) by
# writing them with a Text number format and then clearing the format back off
# so no stray style index is left behind.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") "67.456.62"
$ws.Range("E2").Value = "  -1.26%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.771.84"
$ws.Range("E3").Value = "  +0.17%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  -0.12%  "

# Row 5
Set-TextValue $ws.Range("D5") "593.81"
$ws.Range("E5").Value = "  -0.24%  "

# Row 6
Set-TextValue $ws.Range("D6") "166.45"
$ws.Range("E6").Value = "  -1.13%  "

# Row 7
Set-TextValue $ws.Range("D7") "3.771.10"
$ws.Range("E7").Value = "  +0.21%  "

# Row 8
$ws.Range("E8").Value = "  +0.04%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.521"
$ws.Range("E9").Value = "  -0.29%  "

# Row 10
$ws.Range("E10").Value = "  -0.57%  "

# Row 11
Set-TextValue $ws.Range("D11") "6.40"
$ws.Range("E11").Value = "  -1.41%  "

# Row 12
$ws.Range("E12").Value = "  +0.69%  "

# Row 13
$ws.Range("E13").Value = "  -2.63%  "

# Row 14
Set-TextValue $ws.Range("D14") "35.94"
$ws.Range("E14").Value = "  -1.41%  "

# Row 15
Set-TextValue $ws.Range("D15") "4.406.33"
$ws.Range("E15").Value = "  +0.15%  "

# Row 16
Set-TextValue $ws.Range("D16") "3.791.83"
$ws.Range("E16").Value = "  +0.64%  "

# Row 17
Set-TextValue $ws.Range("D17") "67.560.81"
$ws.Range("E17").Value = "  -1.14%  "

# Row 18
Set-TextValue $ws.Range("D18") "17.79"
$ws.Range("E18").Value = "  -2.14%  "

# Row 19
$ws.Range("E19").Value = "  +0.23%  "

# Row 20
Set-TextValue $ws.Range("D20") "6.94"
$ws.Range("E20").Value = "  -1.62%  "

# Row 21
Set-TextValue $ws.Range("D21") "10.57"
$ws.Range("E21").Value = "  -2.26%  "

# Row 22
Set-TextValue $ws.Range("D22") "460.10"
$ws.Range("E22").Value = "  -1.65%  "

# Row 23
$ws.Range("E23").Value = "  -0.69%  "

# Row 24
Set-TextValue $ws.Range("D24") "0.0000153"
$ws.Range("E24").Value = "  +8.09%  "

# Row 25
Set-TextValue $ws.Range("D25") "83.46"
$ws.Range("E25").Value = "  -0.89%  "

# Row 26
$ws.Range("E26").Value = "  -4.33%  "

# Row 27
Set-TextValue $ws.Range("D27") "11.83"
$ws.Range("E27").Value = "  -2.92%  "

# Row 28
Set-TextValue $ws.Range("D28") "10.04"
$ws.Range("E28").Value = "  -1.61%  "

# Row 29
$ws.Range("E29").Value = "  -0.03%  "

# Row 30
$ws.Range("E30").Value = "  -1.81%  "

# Row 31
Set-TextValue $ws.Range("D31") "7.22"
$ws.Range("E31").Value = "  -2.64%  "

# Row 32
Set-TextValue $ws.Range("D32") "29.71"
$ws.Range("E32").Value = "  -1.04%  "

# Row 33
$ws.Range("E33").Value = "  -2.89%  "

# Row 34
$ws.Range("B34").Value = "Aptos"
$ws.Range("C34").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D34") "9.12"
$ws.Range("E34").Value = "  -1.20%  "

# Row 35
$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws.Range("D35") "1.00"
$ws.Range("E35").Value = "  -0.02%  "

# Row 36
Set-TextValue $ws.Range("D36") "3.724.23"
$ws.Range("E36").Value = "  +0.07%  "

# Row 37
Set-TextValue $ws.Range("D37") "0.1000"
$ws.Range("E37").Value = "  -1.50%  "

# Row 38
$ws.Range("E38").Value = "  -1.70%  "

# Row 39
$ws.Range("E39").Value = "  -0.57%  "

# Row 40
Set-TextValue $ws.Range("D40") "0.995"
$ws.Range("E40").Value = "  -0.46%  "

# Row 41
Set-TextValue $ws.Range("D41") "5.75"
$ws.Range("E41").Value = "  -1.07%  "

# Row 42
$ws.Range("E42").Value = "  -0.18%  "

# Row 43
$ws.Range("E43").Value = "  -0.02%  "

# Row 44
Set-TextValue $ws.Range("D44") "44.22"
$ws.Range("E44").Value = "  +2.11%  "

# Row 45
$ws.Range("E45").Value = "  -2.47%  "

# Row 46
Set-TextValue $ws.Range("D46") "46.91"
$ws.Range("E46").Value = "  +3.32%  "

# Row 47
$ws.Range("E47").Value = "  -3.05%  "

# Row 48
Set-TextValue $ws.Range("D48") "8.34"
$ws.Range("E48").Value = "  -2.61%  "

# Row 49
Set-TextValue $ws.Range("D49") "146.24"
$ws.Range("E49").Value = "  +1.14%  "

# Row 50
Set-TextValue $ws.Range("D50") "389.18"
$ws.Range("E50").Value = "  -3.99%  "

# Row 51
Set-TextValue $ws.Range("D51") "2.756.48"
$ws.Range("E51").Value = "  +3.44%  "
